$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on Hoja1!A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.21 = 8075.31 pesos`n✅ 8075.31 pesos = 2.19 = 977.43 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $newText

# --- Update the rate cells on the "tasas" sheet ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 451.995
$ws2.Range("O10").Value = 3650
$ws2.Range("N12").Value = 3683
